$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells that look like plain decimals (single "." and only
# digits) get written with a leading apostrophe so Excel stores them as plain
# text (matching the source data, which uses "." as a thousands separator for
# larger prices) instead of auto-converting them to numbers; ClearFormats()
# immediately after strips the resulting "number stored as text" styling so the
# cell keeps its original (unstyled) look.

$ws.Range("D2").Value = '26.931.21'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '1.552.39'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  -0.46%  '
$ws.Range("D5").Value = '''206.50'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D8").Value = '''21.95'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("E10").Value = '  +0.56%  '
$ws.Range("D11").Value = '''0.0855'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '1.773.23'
$ws.Range("D13").Value = '1.554.68'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").Value = '''0.519'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("D16").Value = '26.918.90'
$ws.Range("E16").Value = '  -0.37%  '
$ws.Range("D17").Value = '''61.64'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").Value = '0.0₃0712'
$ws.Range("E18").Value = '  +3.25%  '
$ws.Range("D19").Value = '''217.02'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.43%  '
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").Value = '''9.20'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.43%  '
$ws.Range("E24").Value = '  -1.14%  '
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("E30").Value = '  +1.05%  '
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("E33").Value = '  +3.54%  '
$ws.Range("D34").Value = '1.411.59'
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("E35").Value = '  +1.83%  '
$ws.Range("D36").Value = '''0.965'
$ws.Range("D36").ClearFormats()
$ws.Range("E37").Value = '  +0.29%  '
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").Value = '''0.526'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("E40").Value = '  -0.45%  '
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("E42").Value = '  +3.16%  '
$ws.Range("D43").Value = '''2.30'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.23%  '
$ws.Range("D44").Value = '''0.995'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("D45").Value = '''64.46'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.70%  '
$ws.Range("E46").Value = '  -1.48%  '
$ws.Range("D47").Value = '1.687.25'
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("D48").Value = '''87.35'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.27%  '
$ws.Range("D49").Value = '''0.0521'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.74%  '
$ws.Range("E50").Value = '  +2.88%  '
$ws.Range("D51").Value = '''0.0960'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.16%  '
